# Apply property land text cleanup (issue #5: property land done)
# Removes stray spaces / full-width punctuation artifacts introduced by OCR-ish scraping.

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("B2").Value = "臺北市大同區市府段一小段00927000建號"
$wsBuilding.Range("D2").Value = "10000分之1"
$wsBuilding.Range("F2").Value = "98年09月02H"
$wsBuilding.Range("H2").Value = "5000000(為地上權房屋僅有使用權45年無土地所有權無權利持分）"

# 汽車 (Car) sheet
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("B2").Value = "HONDACITY"
$wsCar.Range("E2").Value = "87年10月01日"
$wsCar.Range("G2").Value = "400000(超過五年）"

# 保險 (Insurance) sheet
$wsInsurance = $wb.Worksheets.Item("保險")
$wsInsurance.Range("C3").Value = "富邦人壽真安心醫療養老保險"

# 債務 (Debt) sheet
$wsDebt = $wb.Worksheets.Item("債務")
$wsDebt.Range("D2").Value = "京城銀行忠孝分行臺北市南港區忠孝東路"
$wsDebt.Range("F2").Value = "98年10月27H"
